# Remove the "Device:" label / "${device.deviceName}" value row from the
# report header block, and move the selection cursor onto the now-blank
# A4:B4 cells (commit: "Remove device Name From and change sheet Name").
#
# Rows 5 ("Group:") and 6 ("Period:") below keep their own text untouched;
# clearing A4:B4 simply drops the now-unused "Device:" / "${device.deviceName}"
# shared strings from the workbook on save (sharedStrings count/uniqueCount
# 24 -> 22), which also renumbers every <v> index that follows - Excel does
# this bookkeeping automatically once the cells are cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Blank out the "Device:" label and its "${device.deviceName}" value while
# preserving the existing cell styles (s="1" / s="7").
$ws.Range("A4:B4").ClearContents()

# Move/update the active selection to the cells we just cleared.
$ws.Range("A4:B4").Select()
